$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-12-11 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-12-12 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("58÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "24÷4=", 2) | Out-Null
$d.Content.Find.Execute("72÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "81÷3=", 2) | Out-Null
$d.Content.Find.Execute("25÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "10÷6=", 2) | Out-Null
$d.Content.Find.Execute("76÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2) | Out-Null
$d.Content.Find.Execute("89÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "34÷6=", 2) | Out-Null
$d.Content.Find.Execute("44÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "92÷9=", 2) | Out-Null
$d.Content.Find.Execute("72÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷3=", 2) | Out-Null
$d.Content.Find.Execute("31÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "44÷5=", 2) | Out-Null
$d.Content.Find.Execute("60÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "76÷2=", 2) | Out-Null
$d.Content.Find.Execute("66÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "72÷5=", 2) | Out-Null
$d.Content.Find.Execute("20÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "65÷7=", 2) | Out-Null
$d.Content.Find.Execute("68÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "36÷3=", 2) | Out-Null
$d.Content.Find.Execute("84÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "70÷2=", 2) | Out-Null
$d.Content.Find.Execute("70÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "48÷7=", 2) | Out-Null
$d.Content.Find.Execute("34÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "59÷7=", 2) | Out-Null
$d.Content.Find.Execute("89÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "57÷5=", 2) | Out-Null
$d.Content.Find.Execute("35÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "26÷4=", 2) | Out-Null
$d.Content.Find.Execute("41÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "91÷4=", 2) | Out-Null
$d.Content.Find.Execute("75÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "88÷5=", 2) | Out-Null
$d.Content.Find.Execute("71÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "97÷5=", 2) | Out-Null
$d.Content.Find.Execute("25÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "32÷3=", 2) | Out-Null
$d.Content.Find.Execute("50÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "93÷2=", 2) | Out-Null
$d.Content.Find.Execute("37÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "41÷6=", 2) | Out-Null
$d.Content.Find.Execute("82÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "54÷2=", 2) | Out-Null
$d.Content.Find.Execute("91÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "75÷7=", 2) | Out-Null
